$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record (row) was reported for this market/date and needs to be
# inserted right before the existing row 27, pushing rows 27..63 down to 28..64.
$ws.Rows.Item(27).EntireRow.Insert()

# Fill in the newly inserted row 27 with the new weekly record.
$ws.Range("A27").Value = 5
$ws.Range("B27").Value = "Macroferia Regional de Talca"
$ws.Range("C27").Value = "Maule"
$ws.Range("D27").Value = 44483
$ws.Range("E27").Value = 7
$ws.Range("F27").Value = 100112013
$ws.Range("G27").Value = "Alcachofa"
$ws.Range("H27").Value = "Madrigal"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = 10000
$ws.Range("N27").Value = "$/caja 40 unidades"
$ws.Range("O27").Value = "Provincia del Elquí"
$ws.Range("P27").Value = 250
$ws.Range("Q27").Value = 40
$ws.Range("R27").Value = "Hortaliza"
